$p = $ppt.ActivePresentation
Write-Output $p.Designs.Count
for ($i=1; $i -le $p.Designs.Count; $i++) {
    $d = $p.Designs.Item($i)
    Write-Output "$i : $($d.Name)"
}
